$d = $word.ActiveDocument

# Add the new "Footnote Block Text" paragraph style, based on "Footnote
# Text" (mirrors the existing "Block Text" style, but chained off the
# footnote text style instead of body text).
$wdStyleTypeParagraph = 1
$s = $d.Styles.Add("Footnote Block Text", $wdStyleTypeParagraph)

$s.BaseStyle = "Footnote Text"
$s.NextParagraphStyle = "Footnote Text"
$s.Priority = 9
$s.UnhideWhenUsed = $true
$s.QuickStyle = $true

$pf = $s.ParagraphFormat
$pf.SpaceBefore = 5
$pf.SpaceAfter = 5
$pf.FirstLineIndent = 0
$pf.LeftIndent = 24
$pf.RightIndent = 24
